# Add a new column S (year 2022) to the worksheet, mirroring the existing
# column R (year 2021) formatting, and fill it with the new data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of column R (18) into column S (19) for rows 4-34 so the
# new column inherits the same number formats / borders / fonts as the rest
# of the table, then set the actual values.
$ws.Range("R4:R34").Copy()
$ws.Range("S4:S34").PasteSpecial(-4122)  # xlPasteFormats

# Header value (year)
$ws.Cells.Item(4, 19).Value = 2022

# Data values for rows 5-34 (column S)
$values = @{
    5  = 16.696653653506477
    6  = 17.71894995601205
    7  = 15.612684844888001
    8  = 15.66812062518596
    9  = 16.652881900156387
    10 = 14.667361954014684
    11 = 16.525244796823369
    12 = 19.119250309028729
    13 = 13.749215987119079
    14 = 14.263200620072119
    15 = 16.214093517712168
    16 = 12.189607205170377
    17 = 14.097780631317802
    18 = 16.597474200848456
    19 = 11.56800988291025
    20 = 11.191263248519153
    21 = 11.361761672735106
    22 = 11.015850216858553
    23 = 15.623145704601036
    24 = 18.344423887154832
    25 = 12.847349120106124
    26 = 13.798472231512836
    27 = 15.128863237337196
    28 = 12.394222749619622
    29 = 21.380402934584232
    30 = 19.968977602899539
    31 = 22.891947678227961
    32 = 28.912046224512313
    33 = 31.246721692820181
    34 = 26.427454495987305
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 19).Value = $values[$row]
}

# Update the active cell / selection to T4 (next empty header cell), matching
# the post-edit state recorded in the workbook.
$ws.Range("T4").Select()
